$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on price cells whose new value would otherwise
# be auto-parsed by Excel as a numeric literal (these are plain-text cells
# in the source workbook, t="inlineStr").
$forceTextCells = @("D5","D6","D11","D14","D18","D19","D20","D21","D24","D25","D28","D30","D31","D32","D34","D35","D36","D40","D45","D47","D49","D50","D51")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "57.337.06"
$ws.Range("E2").Value = "  +3.52%  "
$ws.Range("D3").Value = "3.066.86"
$ws.Range("E3").Value = "  +5.58%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "515.00"
$ws.Range("E5").Value = "  +2.70%  "
$ws.Range("D6").Value = "141.05"
$ws.Range("E6").Value = "  +6.70%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +3.56%  "
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("E10").Value = "  +4.27%  "
$ws.Range("D11").Value = "0.371"
$ws.Range("E11").Value = "  +6.91%  "
$ws.Range("D12").Value = "3.593.33"
$ws.Range("E12").Value = "  +5.48%  "
$ws.Range("E13").Value = "  +2.78%  "
$ws.Range("D14").Value = "25.51"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("E15").Value = "  +4.25%  "
$ws.Range("D16").Value = "57.360.20"
$ws.Range("E16").Value = "  +3.55%  "
$ws.Range("D17").Value = "3.070.04"
$ws.Range("E17").Value = "  +5.64%  "
$ws.Range("D18").Value = "5.96"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "12.99"
$ws.Range("E19").Value = "  +4.31%  "
$ws.Range("D20").Value = "8.12"
$ws.Range("E20").Value = "  +6.44%  "
$ws.Range("D21").Value = "337.08"
$ws.Range("E21").Value = "  +7.93%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("E23").Value = "  +3.77%  "
$ws.Range("D24").Value = "65.40"
$ws.Range("E24").Value = "  +4.75%  "
$ws.Range("D25").Value = "0.171"
$ws.Range("E25").Value = "  +7.40%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "0.0₃0941"
$ws.Range("E27").Value = "  +13.16%  "
$ws.Range("D28").Value = "6.42"
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("E29").Value = "  +2.71%  "
$ws.Range("D30").Value = "1.80"
$ws.Range("E30").Value = "  +2.59%  "
$ws.Range("D31").Value = "20.71"
$ws.Range("E31").Value = "  +5.35%  "
$ws.Range("D32").Value = "1.17"
$ws.Range("E32").Value = "  +4.56%  "
$ws.Range("E33").Value = "  +3.42%  "
$ws.Range("D34").Value = "4.52"
$ws.Range("E34").Value = "  +4.09%  "
$ws.Range("D35").Value = "5.86"
$ws.Range("E35").Value = "  +5.66%  "
$ws.Range("D36").Value = "25.94"
$ws.Range("E36").Value = "  +6.81%  "
$ws.Range("E37").Value = "  +5.29%  "
$ws.Range("E38").Value = "  +4.73%  "
$ws.Range("D39").Value = "3.104.94"
$ws.Range("E39").Value = "  +5.68%  "
$ws.Range("D40").Value = "36.95"
$ws.Range("E40").Value = "  +2.16%  "
$ws.Range("E41").Value = "  +5.57%  "
$ws.Range("E42").Value = "  +4.61%  "
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "2.249.47"
$ws.Range("E44").Value = "  +7.65%  "
$ws.Range("D45").Value = "0.0251"
$ws.Range("E45").Value = "  +8.83%  "
$ws.Range("E46").Value = "  +5.05%  "
$ws.Range("D47").Value = "0.948"
$ws.Range("E47").Value = "  +4.87%  "
$ws.Range("E48").Value = "  +8.44%  "
$ws.Range("D49").Value = "5.83"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("D50").Value = "0.0867"
$ws.Range("E50").Value = "  +4.19%  "
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").Value = "0.180"
$ws.Range("E51").Value = "  +5.10%  "
